$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Content fix (commit: "Fixed variables and query errors in Bread from
# TC01 to TC30"): the CasesTab Neo4j query (row 2 / column B) had two
# erroneous trailing lines referencing an undeclared `co`/`cohort` pattern
# that was never MATCHed in this query ("Cohort" column). Remove them so
# the query only returns the columns it actually matches.
# ---------------------------------------------------------------------------
$newCasesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Black and Tan Coonhound']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $newCasesTabQuery

# ---------------------------------------------------------------------------
# View / window cosmetics that accompanied the save: scroll back to the
# top-left of the sheet, select the cell that was just fixed, and zoom in.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 130

# Row heights shrink slightly to match the shorter wrapped text / tighter
# line metrics used when the sheet was resaved.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216
